$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1 ("International Financial Statis"): header row height 16.5 -> 17.25 ---
$ws1.Rows.Item(2).RowHeight = 17.25

# --- Sheet 1: append 12 new monthly rows (Jan 1990 .. Dec 1990) after row 363 ---
# Row 364: Jan 1990
$ws1.Range("A362:F362").Copy($ws1.Range("A364:F364"))
$ws1.Range("A364").Value = "Germany"
$b364 = $ws1.Range("B364")
$b364.NumberFormat = "@"
$b364.Value = "Jan 1990"
$ws1.Range("B362").Copy()
$b364.PasteSpecial(-4122)
$ws1.Range("C364").Value = 78.3694166018768
$ws1.Range("D364").Value = 1.6918
$ws1.Range("E364").Value = 68950.3667233367
$ws1.Range("F364").Value = 66.7174841249911
$ws1.Rows.Item(364).RowHeight = 16.5

# Row 365: Feb 1990
$ws1.Range("A363:F363").Copy($ws1.Range("A365:F365"))
$ws1.Range("A365").Value = "Germany"
$b365 = $ws1.Range("B365")
$b365.NumberFormat = "@"
$b365.Value = "Feb 1990"
$ws1.Range("B363").Copy()
$b365.PasteSpecial(-4122)
$ws1.Range("C365").Value = 78.5858181408675
$ws1.Range("D365").Value = 1.676
$ws1.Range("E365").Value = 69982.7949030696
$ws1.Range("F365").Value = 66.9972272165096
$ws1.Rows.Item(365).RowHeight = 16.5

# Row 366: Mar 1990
$ws1.Range("A362:F362").Copy($ws1.Range("A366:F366"))
$ws1.Range("A366").Value = "Germany"
$b366 = $ws1.Range("B366")
$b366.NumberFormat = "@"
$b366.Value = "Mar 1990"
$ws1.Range("B362").Copy()
$b366.PasteSpecial(-4122)
$ws1.Range("C366").Value = 87.1296714950913
$ws1.Range("D366").Value = 1.7045
$ws1.Range("E366").Value = 69567.6986922397
$ws1.Range("F366").Value = 66.9972272165096
$ws1.Rows.Item(366).RowHeight = 16.5

# Row 367: Apr 1990
$ws1.Range("A363:F363").Copy($ws1.Range("A367:F367"))
$ws1.Range("A367").Value = "Germany"
$b367 = $ws1.Range("B367")
$b367.NumberFormat = "@"
$b367.Value = "Apr 1990"
$ws1.Range("B363").Copy()
$b367.PasteSpecial(-4122)
$ws1.Range("C367").Value = 79.7800340412234
$ws1.Range("D367").Value = 1.6882
$ws1.Range("E367").Value = 69701.9849596987
$ws1.Range("F367").Value = 67.137087460431
$ws1.Rows.Item(367).RowHeight = 16.5

# Row 368: May 1990
$ws1.Range("A362:F362").Copy($ws1.Range("A368:F368"))
$ws1.Range("A368").Value = "Germany"
$b368 = $ws1.Range("B368")
$b368.NumberFormat = "@"
$b368.Value = "May 1990"
$ws1.Range("B362").Copy()
$b368.PasteSpecial(-4122)
$ws1.Range("C368").Value = 81.2708001987146
$ws1.Range("D368").Value = 1.6617
$ws1.Range("E368").Value = 70727.7233721584
$ws1.Range("F368").Value = 67.2769590061902
$ws1.Rows.Item(368).RowHeight = 16.5

# Row 369: Jun 1990
$ws1.Range("A363:F363").Copy($ws1.Range("A369:F369"))
$ws1.Range("A369").Value = "Germany"
$b369 = $ws1.Range("B369")
$b369.NumberFormat = "@"
$b369.Value = "Jun 1990"
$ws1.Range("B363").Copy()
$b369.PasteSpecial(-4122)
$ws1.Range("C369").Value = 81.7597073793972
$ws1.Range("D369").Value = 1.684
$ws1.Range("E369").Value = 71640.0798160129
$ws1.Range("F369").Value = 67.3468947790698
$ws1.Rows.Item(369).RowHeight = 16.5

# Row 370: Jul 1990
$ws1.Range("A362:F362").Copy($ws1.Range("A370:F370"))
$ws1.Range("A370").Value = "Germany"
$b370 = $ws1.Range("B370")
$b370.NumberFormat = "@"
$b370.Value = "Jul 1990"
$ws1.Range("B362").Copy()
$b370.PasteSpecial(-4122)
$ws1.Range("C370").Value = 78.7942048080437
$ws1.Range("D370").Value = 1.6399
$ws1.Range("E370").Value = 71072.396009484
$ws1.Range("F370").Value = 67.3468947790698
$ws1.Rows.Item(370).RowHeight = 16.5

# Row 371: Aug 1990
$ws1.Range("A363:F363").Copy($ws1.Range("A371:F371"))
$ws1.Range("A371").Value = "Germany"
$b371 = $ws1.Range("B371")
$b371.NumberFormat = "@"
$b371.Value = "Aug 1990"
$ws1.Range("B363").Copy()
$b371.PasteSpecial(-4122)
$ws1.Range("C371").Value = 75.9649550575362
$ws1.Range("D371").Value = 1.5707
$ws1.Range("E371").Value = 72425.7385732047
$ws1.Range("F371").Value = 67.5567020977087
$ws1.Rows.Item(371).RowHeight = 16.5

# Row 372: Sep 1990
$ws1.Range("A362:F362").Copy($ws1.Range("A372:F372"))
$ws1.Range("A372").Value = "Germany"
$b372 = $ws1.Range("B372")
$b372.NumberFormat = "@"
$b372.Value = "Sep 1990"
$ws1.Range("B362").Copy()
$b372.PasteSpecial(-4122)
$ws1.Range("C372").Value = 86.4964966217483
$ws1.Range("D372").Value = 1.5697
$ws1.Range("E372").Value = 73197.5736206214
$ws1.Range("F372").Value = 67.7665094163475
$ws1.Rows.Item(372).RowHeight = 16.5

# Row 373: Oct 1990
$ws1.Range("A363:F363").Copy($ws1.Range("A373:F373"))
$ws1.Range("A373").Value = "Germany"
$b373 = $ws1.Range("B373")
$b373.NumberFormat = "@"
$b373.Value = "Oct 1990"
$ws1.Range("B363").Copy()
$b373.PasteSpecial(-4122)
$ws1.Range("C373").Value = 92.643903303446
$ws1.Range("D373").Value = 1.5233
$ws1.Range("E373").Value = 75011.9268302281
$ws1.Range("F373").Value = 68.256048524667
$ws1.Rows.Item(373).RowHeight = 16.5

# Row 374: Nov 1990
$ws1.Range("A362:F362").Copy($ws1.Range("A374:F374"))
$ws1.Range("A374").Value = "Germany"
$b374 = $ws1.Range("B374")
$b374.NumberFormat = "@"
$b374.Value = "Nov 1990"
$ws1.Range("B362").Copy()
$b374.PasteSpecial(-4122)
$ws1.Range("C374").Value = 89.4619991931018
$ws1.Range("D374").Value = 1.487
$ws1.Range("E374").Value = 76166.0747088573
$ws1.Range("F374").Value = 68.1161769789078
$ws1.Rows.Item(374).RowHeight = 16.5

# Row 375: Dec 1990
$ws1.Range("A363:F363").Copy($ws1.Range("A375:F375"))
$ws1.Range("A375").Value = "Germany"
$b375 = $ws1.Range("B375")
$b375.NumberFormat = "@"
$b375.Value = "Dec 1990"
$ws1.Range("B363").Copy()
$b375.PasteSpecial(-4122)
$ws1.Range("C375").Value = 83.8115145639013
$ws1.Range("D375").Value = 1.492
$ws1.Range("E375").Value = 77064.2076742664
$ws1.Range("F375").Value = 68.1861127517874
$ws1.Rows.Item(375).RowHeight = 17.25

# --- Sheet 2 ("Tooltip"): append 12 matching tooltip rows after row 362 ---
# Row 363: Jan 1990
$ws2.Range("A363").Value = "Germany"
$b2_363 = $ws2.Range("B363")
$b2_363.NumberFormat = "@"
$b2_363.Value = "Jan 1990"
$ws2.Range("C363").Value = "Country: Germany`nTime: Jan 1990"
$ws2.Range("D363").Value = "Country: Germany`nTime: Jan 1990"
$ws2.Range("E363").Value = "Country: Germany`nTime: Jan 1990"
$ws2.Range("F363").Value = "Country: Germany`nTime: Jan 1990"
$ws2.Rows.Item(363).AutoFit()

# Row 364: Feb 1990
$ws2.Range("A364").Value = "Germany"
$b2_364 = $ws2.Range("B364")
$b2_364.NumberFormat = "@"
$b2_364.Value = "Feb 1990"
$ws2.Range("C364").Value = "Country: Germany`nTime: Feb 1990"
$ws2.Range("D364").Value = "Country: Germany`nTime: Feb 1990"
$ws2.Range("E364").Value = "Country: Germany`nTime: Feb 1990"
$ws2.Range("F364").Value = "Country: Germany`nTime: Feb 1990"
$ws2.Rows.Item(364).AutoFit()

# Row 365: Mar 1990
$ws2.Range("A365").Value = "Germany"
$b2_365 = $ws2.Range("B365")
$b2_365.NumberFormat = "@"
$b2_365.Value = "Mar 1990"
$ws2.Range("C365").Value = "Country: Germany`nTime: Mar 1990"
$ws2.Range("D365").Value = "Country: Germany`nTime: Mar 1990"
$ws2.Range("E365").Value = "Country: Germany`nTime: Mar 1990"
$ws2.Range("F365").Value = "Country: Germany`nTime: Mar 1990"
$ws2.Rows.Item(365).AutoFit()

# Row 366: Apr 1990
$ws2.Range("A366").Value = "Germany"
$b2_366 = $ws2.Range("B366")
$b2_366.NumberFormat = "@"
$b2_366.Value = "Apr 1990"
$ws2.Range("C366").Value = "Country: Germany`nTime: Apr 1990"
$ws2.Range("D366").Value = "Country: Germany`nTime: Apr 1990"
$ws2.Range("E366").Value = "Country: Germany`nTime: Apr 1990"
$ws2.Range("F366").Value = "Country: Germany`nTime: Apr 1990"
$ws2.Rows.Item(366).AutoFit()

# Row 367: May 1990
$ws2.Range("A367").Value = "Germany"
$b2_367 = $ws2.Range("B367")
$b2_367.NumberFormat = "@"
$b2_367.Value = "May 1990"
$ws2.Range("C367").Value = "Country: Germany`nTime: May 1990"
$ws2.Range("D367").Value = "Country: Germany`nTime: May 1990"
$ws2.Range("E367").Value = "Country: Germany`nTime: May 1990"
$ws2.Range("F367").Value = "Country: Germany`nTime: May 1990"
$ws2.Rows.Item(367).AutoFit()

# Row 368: Jun 1990
$ws2.Range("A368").Value = "Germany"
$b2_368 = $ws2.Range("B368")
$b2_368.NumberFormat = "@"
$b2_368.Value = "Jun 1990"
$ws2.Range("C368").Value = "Country: Germany`nTime: Jun 1990"
$ws2.Range("D368").Value = "Country: Germany`nTime: Jun 1990"
$ws2.Range("E368").Value = "Country: Germany`nTime: Jun 1990"
$ws2.Range("F368").Value = "Country: Germany`nTime: Jun 1990"
$ws2.Rows.Item(368).AutoFit()

# Row 369: Jul 1990
$ws2.Range("A369").Value = "Germany"
$b2_369 = $ws2.Range("B369")
$b2_369.NumberFormat = "@"
$b2_369.Value = "Jul 1990"
$ws2.Range("C369").Value = "Country: Germany`nTime: Jul 1990"
$ws2.Range("D369").Value = "Country: Germany`nTime: Jul 1990"
$ws2.Range("E369").Value = "Country: Germany`nTime: Jul 1990"
$ws2.Range("F369").Value = "Country: Germany`nTime: Jul 1990"
$ws2.Rows.Item(369).AutoFit()

# Row 370: Aug 1990
$ws2.Range("A370").Value = "Germany"
$b2_370 = $ws2.Range("B370")
$b2_370.NumberFormat = "@"
$b2_370.Value = "Aug 1990"
$ws2.Range("C370").Value = "Country: Germany`nTime: Aug 1990"
$ws2.Range("D370").Value = "Country: Germany`nTime: Aug 1990"
$ws2.Range("E370").Value = "Country: Germany`nTime: Aug 1990"
$ws2.Range("F370").Value = "Country: Germany`nTime: Aug 1990"
$ws2.Rows.Item(370).AutoFit()

# Row 371: Sep 1990
$ws2.Range("A371").Value = "Germany"
$b2_371 = $ws2.Range("B371")
$b2_371.NumberFormat = "@"
$b2_371.Value = "Sep 1990"
$ws2.Range("C371").Value = "Country: Germany`nTime: Sep 1990"
$ws2.Range("D371").Value = "Country: Germany`nTime: Sep 1990"
$ws2.Range("E371").Value = "Country: Germany`nTime: Sep 1990"
$ws2.Range("F371").Value = "Country: Germany`nTime: Sep 1990"
$ws2.Rows.Item(371).AutoFit()

# Row 372: Oct 1990
$ws2.Range("A372").Value = "Germany"
$b2_372 = $ws2.Range("B372")
$b2_372.NumberFormat = "@"
$b2_372.Value = "Oct 1990"
$ws2.Range("C372").Value = "Country: Germany`nTime: Oct 1990"
$ws2.Range("D372").Value = "Country: Germany`nTime: Oct 1990"
$ws2.Range("E372").Value = "Country: Germany`nTime: Oct 1990"
$ws2.Range("F372").Value = "Country: Germany`nTime: Oct 1990"
$ws2.Rows.Item(372).AutoFit()

# Row 373: Nov 1990
$ws2.Range("A373").Value = "Germany"
$b2_373 = $ws2.Range("B373")
$b2_373.NumberFormat = "@"
$b2_373.Value = "Nov 1990"
$ws2.Range("C373").Value = "Country: Germany`nTime: Nov 1990"
$ws2.Range("D373").Value = "Country: Germany`nTime: Nov 1990"
$ws2.Range("E373").Value = "Country: Germany`nTime: Nov 1990"
$ws2.Range("F373").Value = "Country: Germany`nTime: Nov 1990"
$ws2.Rows.Item(373).AutoFit()

# Row 374: Dec 1990
$ws2.Range("A374").Value = "Germany"
$b2_374 = $ws2.Range("B374")
$b2_374.NumberFormat = "@"
$b2_374.Value = "Dec 1990"
$ws2.Range("C374").Value = "Country: Germany`nTime: Dec 1990"
$ws2.Range("D374").Value = "Country: Germany`nTime: Dec 1990"
$ws2.Range("E374").Value = "Country: Germany`nTime: Dec 1990"
$ws2.Range("F374").Value = "Country: Germany`nTime: Dec 1990"
$ws2.Rows.Item(374).AutoFit()

# --- Sheet 1: extend merged cell A3:A363 to A3:A375 (done LAST so it does not clear the new cell values) ---
$ws1.Range("A3:A375").Merge()
